# Update countries & provincias Spain
# Re-apply a refreshed data pull from the COVID dashboard (18 Oct 2020, 10:10)
# Armenia and Georgia climbed in the "Casos totales" ranking, pushing the
# countries that used to occupy their rows down by one position each.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "last updated" timestamp -------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 18 de Octubre de 2020 a las 10:10"

# --- Country column re-ranking ---------------------------------------------
# Armenia overtakes Austria and Uzbekistan
$ws.Range("A60").Value = "Armenia"
$ws.Range("A61").Value = "Austria"
$ws.Range("A62").Value = "Uzbekistan"

# Georgia overtakes Madagascar, Albania and Noruega
$ws.Range("A94").Value = "Georgia"
$ws.Range("A95").Value = "Madagascar"
$ws.Range("A96").Value = "Albania"
$ws.Range("A97").Value = "Noruega"

# --- Updated statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ------------------

# Row 4: Estados Unidos
$ws.Range("B4").Value = 8343140
$ws.Range("C4").Value = 475
$ws.Range("D4").Value = 5432452
$ws.Range("E4").Value = 2686405
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 224283

# Row 7: Rusia
$ws.Range("B7").Value = 1399334
$ws.Range("C7").Value = 15099
$ws.Range("D7").Value = 1070576
$ws.Range("E7").Value = 304571
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 185
$ws.Range("H7").Value = 24187

# Row 60: Armenia (new data)
$ws.Range("B60").Value = 64694
$ws.Range("C60").Value = 1694
$ws.Range("D60").Value = 48104
$ws.Range("E60").Value = 15509
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 14
$ws.Range("H60").Value = 1081

# Row 61: Austria (shifted down, keeps its prior data)
$ws.Range("B61").Value = 63134
$ws.Range("C61").Value = 0
$ws.Range("D61").Value = 48771
$ws.Range("E61").Value = 13474
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 889

# Row 62: Uzbekistan (shifted down, keeps its prior data)
$ws.Range("B62").Value = 63124
$ws.Range("C62").Value = 315
$ws.Range("D62").Value = 60080
$ws.Range("E62").Value = 2519
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 3
$ws.Range("H62").Value = 525

# Row 65: Singapur
$ws.Range("B65").Value = 57911
$ws.Range("C65").Value = 7
$ws.Range("D65").Value = 57798
$ws.Range("E65").Value = 85
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 28

# Row 73: Hungria
$ws.Range("B73").Value = 46290
$ws.Range("C73").Value = 1474
$ws.Range("D73").Value = 14088
$ws.Range("E73").Value = 31060
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 33
$ws.Range("H73").Value = 1142

# Row 77: Afganistan
$ws.Range("B77").Value = 40200
$ws.Range("C77").Value = 59
$ws.Range("D77").Value = 33614
$ws.Range("E77").Value = 5094
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 4
$ws.Range("H77").Value = 1492

# Row 94: Georgia (new data)
$ws.Range("B94").Value = 17477
$ws.Range("C94").Value = 1192
$ws.Range("D94").Value = 8060
$ws.Range("E94").Value = 9281
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 8
$ws.Range("H94").Value = 136

# Row 95: Madagascar (shifted down, keeps its prior data)
$ws.Range("B95").Value = 16810
$ws.Range("C95").Value = 0
$ws.Range("D95").Value = 16215
$ws.Range("E95").Value = 357
$ws.Range("F95").Value = 0
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 238

# Row 96: Albania (shifted down, keeps its prior data)
$ws.Range("B96").Value = 16774
$ws.Range("C96").Value = 0
$ws.Range("D96").Value = 10001
$ws.Range("E96").Value = 6325
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 448

# Row 97: Noruega (shifted down, keeps its prior data)
$ws.Range("B97").Value = 16369
$ws.Range("C97").Value = 0
$ws.Range("D97").Value = 11863
$ws.Range("E97").Value = 4228
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 278

# Row 141: Estonia
$ws.Range("B141").Value = 4078
$ws.Range("C141").Value = 26
$ws.Range("D141").Value = 3211
$ws.Range("E141").Value = 799
$ws.Range("F141").Value = 0
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 68
